# Update giáo viên (lecturer) codes in column A from gv1..gv9 to gv01..gv09
# (zero-padded to two digits), leaving gv10..gv19 unchanged, then move the
# active cell selection to H8 (no frozen/scrolled top-left cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

for ($i = 1; $i -le 9; $i++) {
    $ws.Cells.Item($i, 1).Value = "gv0$i"
}

$ws.Activate()
$ws.Range("H8").Select()
